# XRO ENSO forecast update: 2025/05
# - revises the three most recent forecast rows (2025-01, 2025-02, 2025-03)
# - appends a new forecast row for init month 2025-04

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nino34")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

# Row 26 = init 2025-01 (revised lead values)
$row26 = @(-0.757,-0.654,-0.534,-0.49,-0.489,-0.507,-0.518,-0.5,-0.47,-0.467,-0.499,-0.527,-0.513,-0.464,-0.394,-0.305,-0.195,-0.07199999999999999,0.033,0.096)

# Row 27 = init 2025-02 (revised lead values)
$row27 = @(-0.407,-0.349,-0.321,-0.373,-0.453,-0.499,-0.479,-0.436,-0.425,-0.449,-0.468,-0.448,-0.393,-0.327,-0.254,-0.168,-0.074,0.007,0.058,0.09)

# Row 28 = init 2025-03 (revised lead values)
$row28 = @(0.074,0.029,-0.075,-0.198,-0.287,-0.309,-0.301,-0.311,-0.345,-0.373,-0.363,-0.32,-0.267,-0.216,-0.166,-0.116,-0.073,-0.043,-0.024,-0.008)

# Row 29 = init 2025-04 (new forecast row)
$row29 = @(-0.127,-0.169,-0.239,-0.268,-0.239,-0.199,-0.182,-0.178,-0.165,-0.136,-0.105,-0.082,-0.064,-0.041,-0.014,0.011,0.027,0.036,0.048,0.065)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
    $ws.Range($cols[$i] + "28").Value = $row28[$i]
}

# New row 29: label, matching the bold/centered/bordered style used by the
# rest of column A (same look as A2:A28)
$ws.Range("A29").Value = "2025-04"
$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("A29").VerticalAlignment = -4160
$ws.Range("A29").Borders.LineStyle = 1

for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "29")
    $cell.Value = $row29[$i]
    $cell.NumberFormat = $ws.Range($cols[$i] + "28").NumberFormat
}
